# Commit: "Added JB on phy-TGLAB11 to folders"
#
# This adds a new column (N) to Sheet1 of the ComputerFolders workbook for
# user "JB" (username "bothma") working on computer "phy-tglab11" — mirroring
# the existing "phy-tglab11"/hgarcia entries (column G/H) but using bothma's
# own Dropbox-related paths for the rows that differ per-user.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("N1").Value  = "phy-tglab11"
$ws.Range("N2").Value  = "bothma"
$ws.Range("N3").Value  = "Z:\LivemRNA\RawData"
$ws.Range("N4").Value  = "Z:\LivemRNA\FISHAnalysisData"
$ws.Range("N5").Value  = "C:\Users\bothma\Dropbox\LivemRNADatabase"
$ws.Range("N6").Value  = "C:\Users\bothma\Dropbox\LivemRNAData"
$ws.Range("N8").Value  = "C:\Users\bothma\Dropbox\MS2Pausing"
$ws.Range("N9").Value  = "Z:\LivemRNA\mRNADynamics"
$ws.Range("N10").Value = "Z:\FISHDrosophila\Analysis\schnitzcells"

# Leave the window scrolled/split near the new column, with it selected,
# matching where the editor ended up after adding the column.
$excel.ActiveWindow.SplitColumn = 11
$ws.Range("N10").Select()
